# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (want-to-go count) / "最低票价" (min price) figures
# scraped from bilibili show pages, plus a newly-discovered event on the
# "演出" (Performance) sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibition)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Cells.Item(3, 6).Value = 1349
$ws1.Cells.Item(4, 6).Value = 13241
$ws1.Cells.Item(5, 6).Value = 756
$ws1.Cells.Item(10, 6).Value = 1906

# COMICUP GZ06 went from "已售罄" (sold out, text) to back on sale with a
# numeric minimum price.
$ws1.Cells.Item(13, 6).Value = 19415
$ws1.Cells.Item(13, 7).Value = 68

$ws1.Cells.Item(17, 6).Value = 370
$ws1.Cells.Item(20, 6).Value = 155
$ws1.Cells.Item(24, 6).Value = 279
$ws1.Cells.Item(26, 6).Value = 1344
$ws1.Cells.Item(27, 6).Value = 34
$ws1.Cells.Item(28, 6).Value = 366
$ws1.Cells.Item(30, 6).Value = 106

# ---------------------------------------------------------------------
# Sheet "演出" (Performance)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

# Refresh existing figures before the new row shifts everything down.
$ws2.Cells.Item(4, 6).Value = 4463    # LoveLive！电视动画播放十周年纪念巡演
$ws2.Cells.Item(10, 6).Value = 380    # 神山羊2024巡演ENCOUNTER
$ws2.Cells.Item(13, 6).Value = 7      # 《梁祝·卡农》中外经典名曲精选音乐会
$ws2.Cells.Item(15, 6).Value = 17     # 平野宏周粉丝见面会

# A newly-scraped event slots in before "2024-09-01" (the wanuka show),
# pushing every later row down by one.
$ws2.Rows.Item(7).Insert()

$ws2.Cells.Item(7, 1).Value = 6
$ws2.Cells.Item(7, 2).Value = "2024-08-31"
$ws2.Cells.Item(7, 3).Value = "广州·海心沙 •ACGM LIVE SHOW•二次元音乐季"
$ws2.Cells.Item(7, 4).Value = "珠江新城临江大道 广州海心沙亚运公园"
$ws2.Cells.Item(7, 5).Value = "2024.08.31 17:00-09.01 22:00"
$ws2.Cells.Item(7, 6).Value = 0
$ws2.Cells.Item(7, 7).Value = "不可售"
$ws2.Cells.Item(7, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90755"
$ws2.Cells.Item(7, 9).Value = "//i0.hdslb.com/bfs/openplatform/202408/1SrhfYj81723689514036.jpeg"

# Renumber the sequential index column (A) for every row after the insert
# so it keeps counting 0,1,2,... down the sheet (A(n) = n - 2).
for ($r = 8; $r -le 20; $r++) {
    $ws2.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# Sheet "本地生活" (Local Life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Cells.Item(2, 6).Value = 902
$ws3.Cells.Item(3, 6).Value = 4389
$ws3.Cells.Item(4, 6).Value = 67

# ---------------------------------------------------------------------
# Sheet "全部类型" (All Types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Cells.Item(2, 6).Value = 902
$ws4.Cells.Item(5, 6).Value = 1349
$ws4.Cells.Item(6, 6).Value = 13241
$ws4.Cells.Item(8, 6).Value = 756
$ws4.Cells.Item(9, 6).Value = 4389
$ws4.Cells.Item(14, 6).Value = 1906
$ws4.Cells.Item(17, 6).Value = 67

$ws4.Cells.Item(18, 6).Value = 19415
$ws4.Cells.Item(18, 7).Value = 68

$ws4.Cells.Item(20, 6).Value = 4463
$ws4.Cells.Item(29, 6).Value = 380
$ws4.Cells.Item(30, 6).Value = 370
$ws4.Cells.Item(33, 6).Value = 155
$ws4.Cells.Item(38, 6).Value = 7
$ws4.Cells.Item(40, 6).Value = 279
$ws4.Cells.Item(42, 6).Value = 1344
$ws4.Cells.Item(43, 6).Value = 34
$ws4.Cells.Item(44, 6).Value = 17
$ws4.Cells.Item(45, 6).Value = 366
$ws4.Cells.Item(47, 6).Value = 106
